# Apply attendance_reports sync changes to the Session Analysis Results sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: reorder the "Recorded By" list for the ANATOMY session 1 row
$ws.Range("G2").Value = "Veronia.rafat@med.asu.edu.eg, gehanadel@med.asu.edu.eg, System, servinaz@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg"

# Row 3: add System to the "Recorded By" list and update attendance count for ANATOMY session 2
$ws.Range("G3").Value = "Veronia.rafat@med.asu.edu.eg, System, asmaa.reda@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg"
$ws.Range("H3").Value = "103/251"

# Row 10: updated average attendance percentage (kept as literal text, matching source data)
$ws.Range("L10").Value = "'28.1%"

# Row 15: reorder "Recorded By" list and update average attendance percentage for PARASITOLOGY
$ws.Range("G15").Value = "Rania.a.youssef@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg"
$ws.Range("S15").Value = "'28.1%"

# Row 28: reorder "Recorded By" list for PHYSIOLOGY session 1
$ws.Range("G28").Value = "Aya_hamed@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg"
